$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet/tab
$ws.Name = "TASK0-Your responses.downloadlo"

# K column: "Marked" -> "Posted" (shared string used by K2:K4)
$ws.Range("K2").Value = "Posted"
$ws.Range("K3").Value = "Posted"
$ws.Range("K4").Value = "Posted"

# G2: was 1357 with thousands-separator style -> plain number 1.357, default style
$ws.Range("G2").Style = "Normal"
$ws.Range("G2").Value = 1.357

# G3: was text "0.861" -> numeric 0.861
$ws.Range("G3").Value = 0.861

# G4: was text "0.587" -> numeric 0.587
$ws.Range("G4").Value = 0.587
